$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 61.84465033333333
$ws.Range("H2").Value = 185.533951
$ws.Range("I2").Value = 0.03153184209101587
$ws.Range("J2").Value = 0.03153184209101587
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 42.909214
$ws.Range("N2").Value = 128.727642
$ws.Range("O2").Value = 0.2422627718984814
$ws.Range("P2").Value = 0.2422627718984814
$ws.Range("Q2").Value = 2653.705335908171
$ws.Range("R2").Value = 23883.34802317354
$ws.Range("S2").Value = 0.007638991468034714
$ws.Range("T2").Value = 0.007638991468034712

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 61.84465033333333
$ws.Range("H3").Value = 185.533951
$ws.Range("I3").Value = 0.03153184209101587
$ws.Range("J3").Value = 0.03153184209101587
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("N3").Value = 160.0844
$ws.Range("O3").Value = 0.3012755448569878
$ws.Range("P3").Value = 0.3012755448569878
$ws.Range("Q3").Value = 3300.121247273823
$ws.Range("R3").Value = 29701.0912254644
$ws.Range("S3").Value = 0.009499772906315308
$ws.Range("T3").Value = 0.009499772906315306

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 61.84465033333333
$ws.Range("H4").Value = 185.533951
$ws.Range("I4").Value = 0.03153184209101587
$ws.Range("J4").Value = 0.03153184209101587
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 39.98186433333333
$ws.Range("N4").Value = 119.945593
$ws.Range("O4").Value = 0.2257351364921847
$ws.Range("P4").Value = 0.2257351364921847
$ws.Range("Q4").Value = 2472.664419369772
$ws.Range("R4").Value = 22253.97977432795
$ws.Range("S4").Value = 0.007117844678265483
$ws.Range("T4").Value = 0.007117844678265481

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 61.84465033333333
$ws.Range("H5").Value = 185.533951
$ws.Range("I5").Value = 0.03153184209101587
$ws.Range("J5").Value = 0.03153184209101587
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.39145933333334
$ws.Range("N5").Value = 103.174378
$ws.Range("O5").Value = 0.1941720551610951
$ws.Range("P5").Value = 0.1941720551610951
$ws.Range("Q5").Value = 2126.927776923053
$ws.Range("R5").Value = 19142.34999230748
$ws.Range("S5").Value = 0.006122602581827674
$ws.Range("T5").Value = 0.006122602581827672

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 61.84465033333333
$ws.Range("H6").Value = 185.533951
$ws.Range("I6").Value = 0.03153184209101587
$ws.Range("J6").Value = 0.03153184209101587
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.474476
$ws.Range("N6").Value = 19.423428
$ws.Range("O6").Value = 0.03655449159125106
$ws.Range("P6").Value = 0.03655449159125106
$ws.Range("Q6").Value = 400.4117043115587
$ws.Range("R6").Value = 3603.705338804028
$ws.Range("S6").Value = 0.001152630456572696
$ws.Range("T6").Value = 0.001152630456572696

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1361.379069
$ws.Range("H7").Value = 4084.137207
$ws.Range("I7").Value = 0.6941067594101231
$ws.Range("J7").Value = 0.6941067594101232
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.909214
$ws.Range("N7").Value = 128.727642
$ws.Range("O7").Value = 0.2422627718984814
$ws.Range("P7").Value = 0.2422627718984814
$ws.Range("Q7").Value = 58415.70580684176
$ws.Range("R7").Value = 525741.3522615759
$ws.Range("S7").Value = 0.1681562275281688
$ws.Range("T7").Value = 0.1681562275281688

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1361.379069
$ws.Range("H8").Value = 4084.137207
$ws.Range("I8").Value = 0.6941067594101231
$ws.Range("J8").Value = 0.6941067594101232
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("N8").Value = 160.0844
$ws.Range("O8").Value = 0.3012755448569878
$ws.Range("P8").Value = 0.3012755448569878
$ws.Range("Q8").Value = 72645.18381114119
$ws.Range("R8").Value = 653806.6543002708
$ws.Range("S8").Value = 0.209117392130203
$ws.Range("T8").Value = 0.209117392130203

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1361.379069
$ws.Range("H9").Value = 4084.137207
$ws.Range("I9").Value = 0.6941067594101231
$ws.Range("J9").Value = 0.6941067594101232
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 39.98186433333333
$ws.Range("N9").Value = 119.945593
$ws.Range("O9").Value = 0.2257351364921847
$ws.Range("P9").Value = 0.2257351364921847
$ws.Range("Q9").Value = 54430.47324299764
$ws.Range("R9").Value = 489874.2591869787
$ws.Range("S9").Value = 0.1566842840755922
$ws.Range("T9").Value = 0.1566842840755922

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1361.379069
$ws.Range("H10").Value = 4084.137207
$ws.Range("I10").Value = 0.6941067594101231
$ws.Range("J10").Value = 0.6941067594101232
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 34.39145933333334
$ws.Range("N10").Value = 103.174378
$ws.Range("O10").Value = 0.1941720551610951
$ws.Range("P10").Value = 0.1941720551610951
$ws.Range("Q10").Value = 46819.8128887647
$ws.Range("R10").Value = 421378.3159988822
$ws.Range("S10").Value = 0.1347761359758714
$ws.Range("T10").Value = 0.1347761359758714

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1361.379069
$ws.Range("H11").Value = 4084.137207
$ws.Range("I11").Value = 0.6941067594101231
$ws.Range("J11").Value = 0.6941067594101232
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.474476
$ws.Range("N11").Value = 19.423428
$ws.Range("O11").Value = 0.03655449159125106
$ws.Range("P11").Value = 0.03655449159125106
$ws.Range("Q11").Value = 8814.216109142842
$ws.Range("R11").Value = 79327.9449822856
$ws.Range("S11").Value = 0.02537271970028787
$ws.Range("T11").Value = 0.02537271970028787

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 251.007014
$ws.Range("H12").Value = 753.021042
$ws.Range("I12").Value = 0.127977334927537
$ws.Range("J12").Value = 0.1279773349275369
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 42.909214
$ws.Range("N12").Value = 128.727642
$ws.Range("O12").Value = 0.2422627718984814
$ws.Range("P12").Value = 0.2422627718984814
$ws.Range("Q12").Value = 10770.513679227
$ws.Range("R12").Value = 96934.62311304297
$ws.Range("S12").Value = 0.03100414389972544
$ws.Range("T12").Value = 0.03100414389972544

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 251.007014
$ws.Range("H13").Value = 753.021042
$ws.Range("I13").Value = 0.127977334927537
$ws.Range("J13").Value = 0.1279773349275369
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("N13").Value = 160.0844
$ws.Range("O13").Value = 0.3012755448569878
$ws.Range("P13").Value = 0.3012755448569878
$ws.Range("Q13").Value = 13394.10241066053
$ws.Range("R13").Value = 120546.9216959448
$ws.Range("S13").Value = 0.03855644130963891
$ws.Range("T13").Value = 0.0385564413096389

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 251.007014
$ws.Range("H14").Value = 753.021042
$ws.Range("I14").Value = 0.127977334927537
$ws.Range("J14").Value = 0.1279773349275369
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 39.98186433333333
$ws.Range("N14").Value = 119.945593
$ws.Range("O14").Value = 0.2257351364921847
$ws.Range("P14").Value = 0.2257351364921847
$ws.Range("Q14").Value = 10035.7283804631
$ws.Range("R14").Value = 90321.5554241679
$ws.Range("S14").Value = 0.02888898116777359
$ws.Range("T14").Value = 0.02888898116777358

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 251.007014
$ws.Range("H15").Value = 753.021042
$ws.Range("I15").Value = 0.127977334927537
$ws.Range("J15").Value = 0.1279773349275369
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 34.39145933333334
$ws.Range("N15").Value = 103.174378
$ws.Range("O15").Value = 0.1941720551610951
$ws.Range("P15").Value = 0.1941720551610951
$ws.Range("Q15").Value = 8632.497514362432
$ws.Range("R15").Value = 77692.47762926188
$ws.Range("S15").Value = 0.02484962213691965
$ws.Range("T15").Value = 0.02484962213691964

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 251.007014
$ws.Range("H16").Value = 753.021042
$ws.Range("I16").Value = 0.127977334927537
$ws.Range("J16").Value = 0.1279773349275369
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.474476
$ws.Range("N16").Value = 19.423428
$ws.Range("O16").Value = 0.03655449159125106
$ws.Range("P16").Value = 0.03655449159125106
$ws.Range("Q16").Value = 1625.138887974664
$ws.Range("R16").Value = 14626.24999177198
$ws.Range("S16").Value = 0.004678146413479371
$ws.Range("T16").Value = 0.00467814641347937

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 260.0315303333334
$ws.Range("H17").Value = 780.094591
$ws.Range("I17").Value = 0.1325785352324417
$ws.Range("J17").Value = 0.1325785352324417
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 42.909214
$ws.Range("N17").Value = 128.727642
$ws.Range("O17").Value = 0.2422627718984814
$ws.Range("P17").Value = 0.2422627718984814
$ws.Range("Q17").Value = 11157.74858182049
$ws.Range("R17").Value = 100419.7372363844
$ws.Range("S17").Value = 0.03211884343965181
$ws.Range("T17").Value = 0.03211884343965181

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 260.0315303333334
$ws.Range("H18").Value = 780.094591
$ws.Range("I18").Value = 0.1325785352324417
$ws.Range("J18").Value = 0.1325785352324417
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 53.36146666666667
$ws.Range("N18").Value = 160.0844
$ws.Range("O18").Value = 0.3012755448569878
$ws.Range("P18").Value = 0.3012755448569878
$ws.Range("Q18").Value = 13875.66383816449
$ws.Range("R18").Value = 124880.9745434804
$ws.Range("S18").Value = 0.03994267043849523
$ws.Range("T18").Value = 0.03994267043849523

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 260.0315303333334
$ws.Range("H19").Value = 780.094591
$ws.Range("I19").Value = 0.1325785352324417
$ws.Range("J19").Value = 0.1325785352324417
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 39.98186433333333
$ws.Range("N19").Value = 119.945593
$ws.Range("O19").Value = 0.2257351364921847
$ws.Range("P19").Value = 0.2257351364921847
$ws.Range("Q19").Value = 10396.54536817639
$ws.Range("R19").Value = 93568.90831358747
$ws.Range("S19").Value = 0.02992763374662915
$ws.Range("T19").Value = 0.02992763374662914

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 260.0315303333334
$ws.Range("H20").Value = 780.094591
$ws.Range("I20").Value = 0.1325785352324417
$ws.Range("J20").Value = 0.1325785352324417
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 34.39145933333334
$ws.Range("N20").Value = 103.174378
$ws.Range("O20").Value = 0.1941720551610951
$ws.Range("P20").Value = 0.1941720551610951
$ws.Range("Q20").Value = 8942.863800843268
$ws.Range("R20").Value = 80485.7742075894
$ws.Range("S20").Value = 0.02574304665633086
$ws.Range("T20").Value = 0.02574304665633086

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 260.0315303333334
$ws.Range("H21").Value = 780.094591
$ws.Range("I21").Value = 0.1325785352324417
$ws.Range("J21").Value = 0.1325785352324417
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 6.474476
$ws.Range("N21").Value = 19.423428
$ws.Range("O21").Value = 0.03655449159125106
$ws.Range("P21").Value = 0.03655449159125106
$ws.Range("Q21").Value = 1683.567902386439
$ws.Range("R21").Value = 15152.11112147795
$ws.Range("S21").Value = 0.004846340951334673
$ws.Range("T21").Value = 0.004846340951334673

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 27.07732933333333
$ws.Range("H22").Value = 81.231988
$ws.Range("I22").Value = 0.01380552833888228
$ws.Range("J22").Value = 0.01380552833888228
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 42.909214
$ws.Range("N22").Value = 128.727642
$ws.Range("O22").Value = 0.2422627718984814
$ws.Range("P22").Value = 0.2422627718984814
$ws.Range("Q22").Value = 1161.866918912477
$ws.Range("R22").Value = 10456.8022702123
$ws.Range("S22").Value = 0.003344565562900659
$ws.Range("T22").Value = 0.003344565562900659

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 27.07732933333333
$ws.Range("H23").Value = 81.231988
$ws.Range("I23").Value = 0.01380552833888228
$ws.Range("J23").Value = 0.01380552833888228
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 53.36146666666667
$ws.Range("N23").Value = 160.0844
$ws.Range("O23").Value = 0.3012755448569878
$ws.Range("P23").Value = 0.3012755448569878
$ws.Range("Q23").Value = 1444.886006643022
$ws.Range("R23").Value = 13003.9740597872
$ws.Range("S23").Value = 0.004159268072335344
$ws.Range("T23").Value = 0.004159268072335344

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 27.07732933333333
$ws.Range("H24").Value = 81.231988
$ws.Range("I24").Value = 0.01380552833888228
$ws.Range("J24").Value = 0.01380552833888228
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 39.98186433333333
$ws.Range("N24").Value = 119.945593
$ws.Range("O24").Value = 0.2257351364921847
$ws.Range("P24").Value = 0.2257351364921847
$ws.Range("Q24").Value = 1082.602107914321
$ws.Range("R24").Value = 9743.418971228884
$ws.Range("S24").Value = 0.003116392823924315
$ws.Range("T24").Value = 0.003116392823924315

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 27.07732933333333
$ws.Range("H25").Value = 81.231988
$ws.Range("I25").Value = 0.01380552833888228
$ws.Range("J25").Value = 0.01380552833888228
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 34.39145933333334
$ws.Range("N25").Value = 103.174378
$ws.Range("O25").Value = 0.1941720551610951
$ws.Range("P25").Value = 0.1941720551610951
$ws.Range("Q25").Value = 931.2288706226072
$ws.Range("R25").Value = 8381.059835603464
$ws.Range("S25").Value = 0.002680647810145511
$ws.Range("T25").Value = 0.002680647810145511

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 27.07732933333333
$ws.Range("H26").Value = 81.231988
$ws.Range("I26").Value = 0.01380552833888228
$ws.Range("J26").Value = 0.01380552833888228
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 6.474476
$ws.Range("N26").Value = 19.423428
$ws.Range("O26").Value = 0.03655449159125106
$ws.Range("P26").Value = 0.03655449159125106
$ws.Range("Q26").Value = 175.3115189127627
$ws.Range("R26").Value = 1577.803670214864
$ws.Range("S26").Value = 0.0005046540695764505
$ws.Range("T26").Value = 0.0005046540695764505
